# Trade #46 closed at 2026-02-17 13:28:08 - unknown UNKNOWN +0.000%
#
# Updates the Summary / Strategy Status roll-up numbers for the
# MarketMaking strategy and appends the newly closed trade (#46) as a
# new row to both the "All Trades" and "MarketMaking" logs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Summary sheet roll-up figures
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1197.49    # Current Capital
$summary.Range("B4").Value = -2.5       # Total P&L $
$summary.Range("B5").Value = -1.09      # Total P&L %
$summary.Range("B6").Value = 46         # Total Trades
$summary.Range("B8").Value = 25         # Losing Trades
$summary.Range("B9").Value = 39.13      # Win Rate %

# ---------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 97.48999999999999   # Capital
$status.Range("D4").Value = 46                  # Trades
$status.Range("E4").Value = -2.5                # P&L $
$status.Range("F4").Value = -2.51               # P&L %
$status.Range("G4").Value = 39.13               # Win Rate %

# ---------------------------------------------------------------
# 3. Append the new trade row (Trade #46) to a trade-log sheet
# ---------------------------------------------------------------
function Add-TradeRow($ws, $row) {
    $ws.Range("A$row").Value = 46
    # Column B holds a plain "YYYY-MM-DD" string in this workbook (no
    # date number formatting is used anywhere in the sheet) - assigning
    # that literal string via .Value gets auto-parsed into a date
    # serial by Excel, so instead copy it from the identical value
    # already sitting in the row above to keep it as literal text.
    $ws.Range("B" + ($row - 1)).Copy($ws.Range("B$row"))
    $ws.Range("C$row").Value = "13:28:02"
    $ws.Range("D$row").Value = "MarketMaking"
    $ws.Range("E$row").Value = "DOWN"
    $ws.Range("F$row").Value = 0.11
    $ws.Range("G$row").Value = 0.09
    $ws.Range("H$row").Value = "CLOSED"
    $ws.Range("I$row").Value = -18.1818
    $ws.Range("J$row").Value = -0.02
    $ws.Range("K$row").Value = 97.48999999999999
    $ws.Range("L$row").Value = 0
    $ws.Range("M$row").Value = 0
    $ws.Range("N$row").Value = 0.6
    $ws.Range("O$row").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P$row").Value = "early_exit"
    $ws.Range("Q$row").Value = 0.13
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades 47

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking 47
